# Add columns I (header "I0") and J (header "IF") to the sheet, matching
# the style of the existing header row and populating rows 2-68 with data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, thin border) from H1 onto I1/J1,
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-68 for columns I (9) and J (10)
$data = @(
    @(2,8,8),
    @(3,8,8),
    @(4,9,9),
    @(5,8,8),
    @(6,8,8),
    @(7,7,8),
    @(8,8,8),
    @(9,8,8),
    @(10,7,7),
    @(11,8,8),
    @(12,8,8),
    @(13,8,8),
    @(14,7,8),
    @(15,7,7),
    @(16,8,8),
    @(17,8,8),
    @(18,9,9),
    @(19,6,6),
    @(20,8,8),
    @(21,8,8),
    @(22,8,8),
    @(23,6,7),
    @(24,8,8),
    @(25,9,9),
    @(26,8,8),
    @(27,8,8),
    @(28,8,8),
    @(29,7,7),
    @(30,8,8),
    @(31,8,8),
    @(32,6,7),
    @(33,8,8),
    @(34,8,8),
    @(35,8,8),
    @(36,9,9),
    @(37,8,8),
    @(38,8,8),
    @(39,7,8),
    @(40,8,8),
    @(41,9,9),
    @(42,8,8),
    @(43,9,9),
    @(44,8,8),
    @(45,8,8),
    @(46,7,7),
    @(47,8,8),
    @(48,8,8),
    @(49,8,8),
    @(50,8,8),
    @(51,8,8),
    @(52,8,8),
    @(53,8,8),
    @(54,8,8),
    @(55,8,8),
    @(56,8,8),
    @(57,8,8),
    @(58,8,8),
    @(59,8,8),
    @(60,2,2),
    @(61,8,8),
    @(62,5,5),
    @(63,9,9),
    @(64,5,5),
    @(65,5,5),
    @(66,4,4),
    @(67,5,5),
    @(68,5,5)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
